$wb = $excel.ActiveWorkbook

# The "想去人数" (attendance count) figures were refreshed for both the
# "展览" sheet and the "全部类型" sheet. The row numbers differ by one
# between the two sheets because "全部类型" has an extra entry in row 7.

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 819
$ws1.Range("F6").Value = 133
$ws1.Range("F8").Value = 4824
$ws1.Range("F9").Value = 105
$ws1.Range("F10").Value = 5164
$ws1.Range("F12").Value = 1289

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 819
$ws4.Range("F6").Value = 133
$ws4.Range("F9").Value = 4824
$ws4.Range("F10").Value = 105
$ws4.Range("F11").Value = 5164
$ws4.Range("F13").Value = 1289
